$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf = New-Object 'object[,]' 24,5
$bf[0,0] = 1.02
$bf[0,1] = 1.036055746563515
$bf[0,2] = 1.038431795628784
$bf[0,3] = 1.039703526430413
$bf[0,4] = 1.034735289463462
$bf[1,0] = 1.02
$bf[1,1] = 1.037266029392024
$bf[1,2] = 1.039345733560027
$bf[1,3] = 1.040862634453664
$bf[1,4] = 1.03657686658401
$bf[2,0] = 1.02
$bf[2,1] = 1.038048000258962
$bf[2,2] = 1.039936076400179
$bf[2,3] = 1.041611852462715
$bf[2,4] = 1.037767119777883
$bf[3,0] = 1.02
$bf[3,1] = 1.038376466801827
$bf[3,2] = 1.040184011098147
$bf[3,3] = 1.041926635358679
$bf[3,4] = 1.038267185851737
$bf[4,0] = 1.02
$bf[4,1] = 1.038431601839815
$bf[4,2] = 1.040225626108245
$bf[4,3] = 1.041979477854964
$bf[4,4] = 1.038351130932533
$bf[5,0] = 1.02
$bf[5,1] = 1.038052390316176
$bf[5,2] = 1.03993939027691
$bf[5,3] = 1.041616059342664
$bf[5,4] = 1.037773802913903
$bf[6,0] = 1.02
$bf[6,1] = 1.036465010802484
$bf[6,2] = 1.0387408813475
$bf[6,3] = 1.040095420811786
$bf[6,4] = 1.035357949870197
$bf[7,0] = 1.02
$bf[7,1] = 1.03365873835494
$bf[7,2] = 1.036620901784633
$bf[7,3] = 1.037409547529677
$bf[7,4] = 1.031089887679892
$bf[8,0] = 1.02
$bf[8,1] = 1.031781485749073
$bf[8,2] = 1.035201998931867
$bf[8,3] = 1.035614481718923
$bf[8,4] = 1.028236367504404
$bf[9,0] = 1.02
$bf[9,1] = 1.030967032241109
$bf[9,2] = 1.034586235758235
$bf[9,3] = 1.0348360810206
$bf[9,4] = 1.026998665664421
$bf[10,0] = 1.02
$bf[10,1] = 1.030664263395028
$bf[10,2] = 1.034357304838125
$bf[10,3] = 1.034546774526069
$bf[10,4] = 1.026538597082158
$bf[11,0] = 1.02
$bf[11,1] = 1.030729219496402
$bf[11,2] = 1.034406420813047
$bf[11,3] = 1.034608839693432
$bf[11,4] = 1.026637298562766
$bf[12,0] = 1.02
$bf[12,1] = 1.030942010294896
$bf[12,2] = 1.034567316541543
$bf[12,3] = 1.034812170428921
$bf[12,4] = 1.026960643086902
$bf[13,0] = 1.02
$bf[13,1] = 1.031073085156569
$bf[13,2] = 1.034666421905155
$bf[13,3] = 1.034937426030732
$bf[13,4] = 1.027159822031367
$bf[14,0] = 1.02
$bf[14,1] = 1.031835504357514
$bf[14,2] = 1.035242835975483
$bf[14,3] = 1.035666117495668
$bf[14,4] = 1.028318464061173
$bf[15,0] = 1.02
$bf[15,1] = 1.032313319981943
$bf[15,2] = 1.035604036414877
$bf[15,3] = 1.036122901772965
$bf[15,4] = 1.02904467503433
$bf[16,0] = 1.02
$bf[16,1] = 1.032591868646771
$bf[16,2] = 1.035814586676061
$bf[16,3] = 1.036389228056684
$bf[16,4] = 1.029468059106867
$bf[17,0] = 1.02
$bf[17,1] = 1.032686820741037
$bf[17,2] = 1.035886356613042
$bf[17,3] = 1.036480020217219
$bf[17,4] = 1.029612388214067
$bf[18,0] = 1.02
$bf[18,1] = 1.032262070765459
$bf[18,2] = 1.035565296689011
$bf[18,3] = 1.036073904376147
$bf[18,4] = 1.028966780520638
$bf[19,0] = 1.02
$bf[19,1] = 1.030879355504044
$bf[19,2] = 1.034519942558316
$bf[19,3] = 1.034752299431716
$bf[19,4] = 1.026865435435733
$bf[20,0] = 1.02
$bf[20,1] = 1.030008569865266
$bf[20,2] = 1.033861475656042
$bf[20,3] = 1.033920345917816
$bf[20,4] = 1.025542314144583
$bf[21,0] = 1.02
$bf[21,1] = 1.030470326054173
$bf[21,2] = 1.034210657378568
$bf[21,3] = 1.034361477360753
$bf[21,4] = 1.026243912748048
$bf[22,0] = 1.02
$bf[22,1] = 1.032285228553088
$bf[22,2] = 1.035582801910791
$bf[22,3] = 1.036096044524906
$bf[22,4] = 1.029001978322787
$bf[23,0] = 1.02
$bf[23,1] = 1.034385335095648
$bf[23,2] = 1.037169938938617
$bf[23,3] = 1.038104682810893
$bf[23,4] = 1.032194665472047
$ws.Range("B2:F25").Value = $bf

$inn = New-Object 'object[,]' 24,6
$inn[0,0] = 1.038921210627312
$inn[0,1] = 1.041166149162644
$inn[0,2] = 1.041219846780695
$inn[0,3] = 1.042487960110096
$inn[0,4] = 1.03753391026133
$inn[0,5] = 1.017538200654688
$inn[1,0] = 1.039315242192911
$inn[1,1] = 1.042018923185883
$inn[1,2] = 1.041943522176788
$inn[1,3] = 1.043456422144053
$inn[1,4] = 1.039181990904166
$inn[1,5] = 1.017832929549581
$inn[2,0] = 1.039568257392115
$inn[2,1] = 1.042569068801202
$inn[2,2] = 1.042410127747965
$inn[2,3] = 1.044081707712323
$inn[2,4] = 1.040246623809783
$inn[2,5] = 1.01802279790226
$inn[3,0] = 1.03967415996568
$inn[3,1] = 1.042799956822316
$inn[3,2] = 1.042605893760026
$inn[3,3] = 1.044344252274389
$inn[3,4] = 1.04069377847455
$inn[3,5] = 1.018102418245239
$inn[4,0] = 1.039691914254596
$inn[4,1] = 1.042838700975372
$inn[4,2] = 1.042638740649646
$inn[4,3] = 1.04438831567529
$inn[4,4] = 1.040768833558941
$inn[4,5] = 1.018115775134038
$inn[5,0] = 1.039569674292573
$inn[5,1] = 1.042572155479158
$inn[5,2] = 1.042412745131726
$inn[5,3] = 1.044085217120277
$inn[5,4] = 1.040252600337975
$inn[5,5] = 1.01802386257902
$inn[6,0] = 1.039054780147973
$inn[6,1] = 1.04145469305844
$inn[6,2] = 1.041464761737402
$inn[6,3] = 1.042815543011482
$inn[6,4] = 1.03809126282405
$inn[6,5] = 1.017637980400865
$inn[7,0] = 1.038132462015237
$inn[7,1] = 1.039472748088476
$inn[7,2] = 1.039781455430344
$inn[7,3] = 1.040567534027827
$inn[7,4] = 1.034268566513044
$inn[7,5] = 1.016951520279107
$inn[8,0] = 1.037507382198087
$inn[8,1] = 1.038142612440323
$inn[8,2] = 1.038650451094997
$inn[8,3] = 1.039061458270402
$inn[8,4] = 1.031709935899169
$inn[8,5] = 1.016489454100309
$inn[9,0] = 1.037234271205336
$inn[9,1] = 1.037564504954822
$inn[9,2] = 1.038158588297747
$inn[9,3] = 1.038407504787177
$inn[9,4] = 1.03059946259778
$inn[9,5] = 1.016288309748808
$inn[10,0] = 1.037132455716783
$inn[10,1] = 1.037349442894079
$inn[10,2] = 1.037975565098822
$inn[10,3] = 1.038164320586858
$inn[10,4] = 1.030186584655594
$inn[10,5] = 1.01621343421282
$inn[11,0] = 1.037154312259651
$inn[11,1] = 1.03739558931035
$inn[11,2] = 1.038014838887922
$inn[11,3] = 1.038216497018589
$inn[11,4] = 1.03027516662334
$inn[11,5] = 1.016229502602014
$inn[12,0] = 1.037225862664198
$inn[12,1] = 1.037546734556958
$inn[12,2] = 1.038143466161089
$inn[12,3] = 1.038387408776806
$inn[12,4] = 1.030565342185746
$inn[12,5] = 1.01628212381724
$inn[13,0] = 1.0372698981469
$inn[13,1] = 1.03763981664857
$inn[13,2] = 1.038222674706663
$inn[13,3] = 1.038492676345782
$inn[13,4] = 1.030744075711644
$inn[13,5] = 1.016314524035759
$inn[14,0] = 1.037525455923778
$inn[14,1] = 1.038180933875827
$inn[14,2] = 1.038683049234992
$inn[14,3] = 1.039104820440967
$inn[14,4] = 1.031783579140043
$inn[14,5] = 1.016502780803374
$inn[15,0] = 1.037685103787006
$inn[15,1] = 1.038519783859959
$inn[15,2] = 1.038971257254207
$inn[15,3] = 1.039488313831489
$inn[15,4] = 1.032434935354669
$inn[15,5] = 1.016620582802849
$inn[16,0] = 1.037777987706779
$inn[16,1] = 1.038717222167368
$inn[16,2] = 1.039139158788092
$inn[16,3] = 1.039711824435726
$inn[16,4] = 1.032814613371512
$inn[16,5] = 1.016689191889986
$inn[17,0] = 1.037809618763986
$inn[17,1] = 1.038784508516734
$inn[17,2] = 1.039196374146524
$inn[17,3] = 1.039788006292843
$inn[17,4] = 1.032944032243139
$inn[17,5] = 1.016712568414101
$inn[18,0] = 1.037667999507206
$inn[18,1] = 1.038483449922003
$inn[18,2] = 1.038940356524718
$inn[18,3] = 1.039447186713101
$inn[18,4] = 1.032365076644264
$inn[18,5] = 1.016607954408876
$inn[19,0] = 1.037204803082583
$inn[19,1] = 1.037502235115205
$inn[19,2] = 1.038105597591348
$inn[19,3] = 1.038337087199442
$inn[19,4] = 1.030479903823994
$inn[19,5] = 1.016266632651602
$inn[20,0] = 1.036911431829859
$inn[20,1] = 1.036883411238966
$inn[20,2] = 1.037578878406593
$inn[20,3] = 1.037637520116294
$inn[20,4] = 1.029292307649869
$inn[20,5] = 1.016051094673023
$inn[21,0] = 1.037067157197925
$inn[21,1] = 1.037211642534204
$inn[21,2] = 1.03785828089752
$inn[21,3] = 1.038008527336496
$inn[21,4] = 1.029922098381577
$inn[21,5] = 1.016165444545223
$inn[22,0] = 1.037675728923986
$inn[22,1] = 1.038499868304377
$inn[22,2] = 1.038954319868478
$inn[22,3] = 1.039465770826805
$inn[22,4] = 1.03239664355058
$inn[22,5] = 1.01661366095417
$inn[23,0] = 1.038372693243779
$inn[23,1] = 1.039986670758239
$inn[23,2] = 1.040218169108321
$inn[23,3] = 1.041149986703071
$inn[23,4] = 1.035258568572324
$inn[23,5] = 1.017129762002455
$ws.Range("I2:N25").Value = $inn
